# Update column G ("K") values for rows 2-12 on the active sheet.
# These reflect regenerated save_data values (K instead of Strike#).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValues = @{
    2  = 7
    3  = 5
    4  = 2
    5  = 2
    6  = 6
    7  = 1
    8  = 3
    9  = 3
    10 = 4
    11 = 3
    12 = 2
}

foreach ($row in $newValues.Keys) {
    $ws.Range("G$row").Value = $newValues[$row]
}
